$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.975.94"
$ws.Range("E2").Value = "  -0.19%  "
$ws.Range("D3").Value = "1.676.51"
$ws.Range("E3").Value = "  +0.13%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "215.19"
$ws.Range("E5").Value = "  -0.42%  "
$ws.Range("E6").Value = "  +1.43%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("E8").Value = "  -0.14%  "
$ws.Range("E9").Value = "  +0.34%  "
$ws.Range("E10").Value = "  +0.77%  "
$ws.Range("E11").Value = "  +0.03%  "
$ws.Range("D12").Value = "1.914.22"
$ws.Range("E12").Value = "  +0.21%  "
$ws.Range("D13").Value = "1.737.46"
$ws.Range("E13").Value = "  +3.76%  "
$ws.Range("E14").Value = "  -0.10%  "
$ws.Range("E15").Value = "  +1.27%  "
$ws.Range("E16").Value = "  +0.01%  "
$ws.Range("D17").Value = "26.992.97"
$ws.Range("E17").Value = "  -0.21%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "236.94"
$ws.Range("E18").Value = "  +0.43%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "8.06"
$ws.Range("E19").Value = "  +4.38%  "
$ws.Range("E20").Value = "  -0.62%  "
$ws.Range("E21").Value = "  +0.01%  "
$ws.Range("E22").Value = "  -0.67%  "
$ws.Range("E23").Value = "  -0.96%  "
$ws.Range("E24").Value = "  -1.48%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "145.64"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.23"
$ws.Range("E26").Value = "  +1.00%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "16.08"
$ws.Range("E27").Value = "  +0.81%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.112"
$ws.Range("E28").Value = "  -1.22%  "
$ws.Range("E29").Value = "  -0.09%  "
$ws.Range("E31").Value = "  -0.50%  "
$ws.Range("E32").Value = "  +0.20%  "
$ws.Range("D33").Value = "1.485.02"
$ws.Range("E33").Value = "  +1.02%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.16"
$ws.Range("E34").Value = "  +0.78%  "
$ws.Range("E35").Value = "  +4.75%  "
$ws.Range("E36").Value = "  +0.17%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.585"
$ws.Range("E37").Value = "  +1.88%  "
$ws.Range("E38").Value = "  +2.95%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.902"
$ws.Range("E39").Value = "  +0.56%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.87"
$ws.Range("E40").Value = "  -4.16%  "
$ws.Range("E41").Value = "  +1.19%  "
$ws.Range("E42").Value = "  +0.04%  "
$ws.Range("E43").Value = "  +1.96%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "67.45"
$ws.Range("E44").Value = "  +1.52%  "
$ws.Range("D45").Value = "1.819.80"
$ws.Range("E45").Value = "  +0.00%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.782"
$ws.Range("E46").Value = "  +0.34%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "90.65"
$ws.Range("E47").Value = "  +0.46%  "
$ws.Range("E48").Value = "  +2.02%  "
$ws.Range("E49").Value = "  -0.43%  "
$ws.Range("E50").Value = "  +1.79%  "
$ws.Range("E51").Value = "  +0.40%  "
